# Update automatic: dades i banners [2026-02-15 20:50]
# Refresh each weather-station row with the newly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (percentage humidity readings) - force them to stay plain text first so
# the stored value matches the scraped string exactly ("71%", not 0.71).
$textCells = @("H10","H11","H12","H13","H16","H17","H19","H20","H21","H23","H30","H32","H33","H36","H38","H39","H40")
foreach ($tc in $textCells) {
    $ws.Range($tc).NumberFormat = "@"
}

# Row 2
$ws.Range("E2").Value = '2026-02-15 20:48:39'

# Row 3
$ws.Range("E3").Value = '2026-02-15 20:48:41'
$ws.Range("G3").Value = '199 cm'
$ws.Range("I3").Value = '2.2 mm'
$ws.Range("O3").Value = '-5.0 °C'

# Row 4
$ws.Range("E4").Value = '2026-02-15 20:48:44'
$ws.Range("J4").Value = '1015.8 hPa'

# Row 5
$ws.Range("E5").Value = '2026-02-15 20:48:46'
$ws.Range("I5").Value = '6.7 mm'
$ws.Range("L5").Value = '43.9 km/h - 330º 20:04 TU'
$ws.Range("O5").Value = '-4.5 °C'

# Row 6
$ws.Range("E6").Value = '2026-02-15 20:48:49'

# Row 7
$ws.Range("E7").Value = '2026-02-15 20:48:51'
$ws.Range("J7").Value = '1015.9 hPa'
$ws.Range("O7").Value = '11.8 °C'

# Row 8
$ws.Range("E8").Value = '2026-02-15 20:48:54'

# Row 9
$ws.Range("E9").Value = '2026-02-15 20:48:57'

# Row 10
$ws.Range("E10").Value = '2026-02-15 20:48:59'
$ws.Range("H10").Value = '71%'
$ws.Range("O10").Value = '7.4 °C'

# Row 11
$ws.Range("E11").Value = '2026-02-15 20:49:02'
$ws.Range("H11").Value = '44%'
$ws.Range("O11").Value = '7.2 °C'

# Row 12
$ws.Range("E12").Value = '2026-02-15 20:49:04'
$ws.Range("H12").Value = '57%'
$ws.Range("O12").Value = '10.8 °C'

# Row 13
$ws.Range("E13").Value = '2026-02-15 20:49:06'
$ws.Range("H13").Value = '37%'

# Row 14
$ws.Range("E14").Value = '2026-02-15 20:49:09'

# Row 15
$ws.Range("E15").Value = '2026-02-15 20:49:11'

# Row 16
$ws.Range("E16").Value = '2026-02-15 20:49:14'
$ws.Range("H16").Value = '62%'
$ws.Range("I16").Value = '1.0 mm'

# Row 17
$ws.Range("E17").Value = '2026-02-15 20:49:17'
$ws.Range("H17").Value = '39%'

# Row 18
$ws.Range("E18").Value = '2026-02-15 20:49:19'
$ws.Range("J18").Value = '1016.0 hPa'

# Row 19
$ws.Range("E19").Value = '2026-02-15 20:49:22'
$ws.Range("H19").Value = '73%'

# Row 20
$ws.Range("E20").Value = '2026-02-15 20:49:24'
$ws.Range("H20").Value = '61%'
$ws.Range("L20").Value = '74.2 km/h - 333º 20:28 TU'

# Row 21
$ws.Range("E21").Value = '2026-02-15 20:49:27'
$ws.Range("H21").Value = '39%'

# Row 22
$ws.Range("E22").Value = '2026-02-15 20:49:29'

# Row 23
$ws.Range("E23").Value = '2026-02-15 20:49:32'
$ws.Range("H23").Value = '66%'
$ws.Range("I23").Value = '3.8 mm'

# Row 24
$ws.Range("E24").Value = '2026-02-15 20:49:34'

# Row 25
$ws.Range("E25").Value = '2026-02-15 20:49:37'
$ws.Range("O25").Value = '-1.5 °C'

# Row 26
$ws.Range("E26").Value = '2026-02-15 20:49:39'

# Row 27
$ws.Range("E27").Value = '2026-02-15 20:49:41'

# Row 28
$ws.Range("E28").Value = '2026-02-15 20:49:44'

# Row 29
$ws.Range("E29").Value = '2026-02-15 20:49:47'

# Row 30
$ws.Range("E30").Value = '2026-02-15 20:49:49'
$ws.Range("H30").Value = '56%'
$ws.Range("J30").Value = '1015.3 hPa'
$ws.Range("O30").Value = '9.8 °C'

# Row 31
$ws.Range("E31").Value = '2026-02-15 20:49:52'

# Row 32
$ws.Range("E32").Value = '2026-02-15 20:49:54'
$ws.Range("H32").Value = '82%'
$ws.Range("L32").Value = '49.0 km/h - 268º 20:15 TU'
$ws.Range("O32").Value = '3.8 °C'

# Row 33
$ws.Range("E33").Value = '2026-02-15 20:49:57'
$ws.Range("H33").Value = '42%'
$ws.Range("O33").Value = '5.9 °C'

# Row 34
$ws.Range("E34").Value = '2026-02-15 20:49:59'
$ws.Range("M34").Value = '4.7 °C 20:12 TU'
$ws.Range("O34").Value = '1.2 °C'

# Row 35
$ws.Range("E35").Value = '2026-02-15 20:50:02'
$ws.Range("O35").Value = '4.2 °C'

# Row 36
$ws.Range("E36").Value = '2026-02-15 20:50:04'
$ws.Range("H36").Value = '49%'

# Row 37
$ws.Range("E37").Value = '2026-02-15 20:50:07'
$ws.Range("J37").Value = '1016.4 hPa'

# Row 38
$ws.Range("E38").Value = '2026-02-15 20:50:09'
$ws.Range("H38").Value = '68%'

# Row 39
$ws.Range("E39").Value = '2026-02-15 20:50:12'
$ws.Range("H39").Value = '59%'
$ws.Range("O39").Value = '-2.7 °C'

# Row 40
$ws.Range("E40").Value = '2026-02-15 20:50:14'
$ws.Range("H40").Value = '38%'
$ws.Range("O40").Value = '8.8 °C'

# Row 41
$ws.Range("E41").Value = '2026-02-15 20:50:17'

# Row 42
$ws.Range("E42").Value = '2026-02-15 20:50:19'

# Row 43
$ws.Range("E43").Value = '2026-02-15 20:50:22'
$ws.Range("O43").Value = '6.4 °C'

# Row 44
$ws.Range("E44").Value = '2026-02-15 20:50:24'
$ws.Range("I44").Value = '4.4 mm'
$ws.Range("M44").Value = '-0.8 °C 20:02 TU'
$ws.Range("O44").Value = '-3.9 °C'

# Row 45
$ws.Range("E45").Value = '2026-02-15 20:50:27'
$ws.Range("I45").Value = '2.9 mm'
$ws.Range("J45").Value = '1023.3 hPa'
$ws.Range("O45").Value = '1.1 °C'

# Row 46
$ws.Range("E46").Value = '2026-02-15 20:50:29'
$ws.Range("O46").Value = '11.8 °C'
